$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The underlying observation records for two pairs of rows got reordered:
#   row 63 <-> row 66
#   row 64 <-> row 65
# Swap the per-record fields (the columns that actually differ between the
# two rows of each pair) while leaving shared/unaffected columns (dates,
# municipality, reporter, etc.) untouched.

function Swap-Cell($ws, $r1, $r2, $col) {
    $cell1 = $ws.Range($col + $r1)
    $cell2 = $ws.Range($col + $r2)

    $v1 = $cell1.Value2
    $v2 = $cell2.Value2

    $cell1.Value = $v2
    $cell2.Value = $v1
}

$cols6366 = @("A","B","E","F","G","H","I","J","P","Q","R","S","Z","AB")
foreach ($col in $cols6366) {
    Swap-Cell $ws 63 66 $col
}

$cols6465 = @("A","B","E","F","G","H","Z","AB")
foreach ($col in $cols6465) {
    Swap-Cell $ws 64 65 $col
}

# Column L: the empty placeholder cell present on row 65 moves to row 64.
$ws.Range("L64").Value = ""
$ws.Range("L65").Value = ""

# Column I ("Antal") is always stored as text in this sheet, even when the
# text is a bare number (e.g. "3") -- force it back to text so it doesn't
# get reinterpreted as a real number.
$ws.Range("I63").Value = "'3"
